$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds plain text values (e.g. "66.723.26", "1.00").
# Force text format per-cell first so Excel does not auto-convert numeric-
# looking strings ("1.00" -> 1, "142.67" -> 142.67 as a number, etc.)
# and thereby lose formatting/precision.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.723.26'
$ws.Range('E2').Value = '  -0.24%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.786.55'
$ws.Range('E3').Value = '  -2.45%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '436.60'
$ws.Range('E5').Value = '  +1.68%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.67'
$ws.Range('E6').Value = '  +8.09%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.620'
$ws.Range('E7').Value = '  +1.13%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.736'
$ws.Range('E9').Value = '  +0.58%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.153'
$ws.Range('E10').Value = '  -8.49%  '
$ws.Range('E11').Value = '  -13.25%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '43.06'
$ws.Range('E12').Value = '  +5.16%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.44'
$ws.Range('E13').Value = '  +3.59%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.380.36'
$ws.Range('E14').Value = '  -2.67%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.84'
$ws.Range('E15').Value = '  -5.53%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.788.78'
$ws.Range('E16').Value = '  -2.30%  '
$ws.Range('E17').Value = '  -0.67%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '19.91'
$ws.Range('E18').Value = '  +1.18%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.14'
$ws.Range('E19').Value = '  +6.03%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '66.745.19'
$ws.Range('E20').Value = '  -0.89%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '417.98'
$ws.Range('E21').Value = '  +2.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '14.53'
$ws.Range('E22').Value = '  +0.38%  '
$ws.Range('E23').Value = '  +7.64%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '86.04'
$ws.Range('E24').Value = '  +0.69%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '37.24'
$ws.Range('E25').Value = '  -2.42%  '
$ws.Range('E26').Value = '  +4.75%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '5.62'
$ws.Range('E27').Value = '  -0.92%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.59'
$ws.Range('E28').Value = '  +33.79%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.81'
$ws.Range('E29').Value = '  +2.29%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '720.94'
$ws.Range('E30').Value = '  +4.42%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '13.79'
$ws.Range('E31').Value = '  +10.64%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.134'
$ws.Range('E32').Value = '  +9.54%  '
$ws.Range('E33').Value = '  -0.26%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '43.84'
$ws.Range('E34').Value = '  +13.17%  '
$ws.Range('E35').Value = '  +1.78%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  +0.00%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.62'
$ws.Range('E37').Value = '  +24.23%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '56.48'
$ws.Range('E38').Value = '  +2.11%  '
$ws.Range('E39').Value = '  +4.34%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.65'
$ws.Range('E40').Value = '  +37.04%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.91'
$ws.Range('E41').Value = '  -5.16%  '
$ws.Range('E42').Value = '  +3.02%  '
$ws.Range('B43').Value = 'PEPE'
$ws.Range('C43').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0₃0675'
$ws.Range('E43').Value = '  -16.43%  '
$ws.Range('B44').Value = 'ApeXProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.34'
$ws.Range('E44').Value = '  +7.69%  '
$ws.Range('E45').Value = '  -0.19%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.327'
$ws.Range('E46').Value = '  +11.77%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.31'
$ws.Range('E47').Value = '  +1.07%  '
$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.08'
$ws.Range('E48').Value = '  -0.37%  '
$ws.Range('B49').Value = 'WEMIXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.64'
$ws.Range('E49').Value = '  +3.93%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '142.57'
$ws.Range('E50').Value = '  -3.63%  '
$ws.Range('E51').Value = '  +2.01%  '
